$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Beads": replace the single example bead-file row with three rows
# (BMIN / BMAX calibration entries).
# ---------------------------------------------------------------------------
$beads = $wb.Worksheets.Item("Beads")

$beads.Range("A2:G4").ClearContents()

$beads.Range("A2").Value = "B0001"
$beads.Range("B2").Value = "FC001"
$beads.Range("C2").Value = "./FCFiles/sample001.fcs"
$beads.Range("D2").Value = "AJ01"
$beads.Range("E2").Value = "0, 792, 2079, 6588, 16471, 47497, 137049, 271647"
$beads.Range("F2").Value = 0.85
$beads.Range("G2").Value = "FL1, FL3"

$beads.Range("A3").Value = "BMIN"
$beads.Range("B3").Value = "FC001"
$beads.Range("C3").Value = "./FCFiles/min/sample001.fcs"
$beads.Range("D3").Value = "AK02"
$beads.Range("E3").Value = "0, 771, 2106, 6262, 15183, 45292, 136258, 291042"
$beads.Range("F3").Value = 0.85
$beads.Range("G3").Value = "FL1, FL3"

$beads.Range("A4").Value = "BMAX"
$beads.Range("B4").Value = "FC001"
$beads.Range("C4").Value = "./FCFiles/max/sample002.fcs"
$beads.Range("D4").Value = "AJ01"
$beads.Range("E4").Value = "0, 792, 2079, 6588, 16471, 47497, 137049, 271647"
$beads.Range("F4").Value = 0.85
$beads.Range("G4").Value = "FL1, FL3"

# ---------------------------------------------------------------------------
# Sheet "Samples": header rename (Plasmid -> Plasmids, IPTG (mM) -> DAPG
# (uM)) and a fully new set of 12 data rows.
# ---------------------------------------------------------------------------
$samples = $wb.Worksheets.Item("Samples")

$samples.Range("H1").Value = "Plasmids"
$samples.Range("I1").Value = "DAPG (uM)"

$samples.Range("A2:I13").ClearContents()

$sampleRows = @(
    @("S0001", "FC001", "B0001", "./FCFiles/sample006.fcs", "MEF", 0.85, "sJS1123", "pJS0355, pJS0304, pSC31_3", 0),
    @("S0002", "FC001", "B0001", "./FCFiles/sample007.fcs", "MEF", 0.85, "sJS1123", "pJS0355, pJS0304, pSC31_3", 2.332362),
    @("S0003", "FC001", "B0001", "./FCFiles/sample008.fcs", "MEF", 0.85, "sJS1123", "pJS0355, pJS0304, pSC31_3", 4.363449),
    @("S0004", "FC001", "B0001", "./FCFiles/sample009.fcs", "MEF", 0.85, "sJS1123", "pJS0355, pJS0304, pSC31_3", 8.163265),
    @("S0005", "FC001", "B0001", "./FCFiles/sample010.fcs", "MEF", 0.85, "sJS1123", "pJS0355, pJS0304, pSC31_3", 15.27207),
    @("S0006", "FC001", "B0001", "./FCFiles/sample011.fcs", "MEF", 0.85, "sJS1123", "pJS0355, pJS0304, pSC31_3", 28.57143),
    @("S0007", "FC001", "B0001", "./FCFiles/sample012.fcs", "MEF", 0.85, "sJS1123", "pJS0355, pJS0304, pSC31_3", 53.45225),
    @("S0008", "FC001", "B0001", "./FCFiles/sample013.fcs", "MEF", 0.85, "sJS1123", "pJS0355, pJS0304, pSC31_3", 100),
    @("S0009", "FC001", "B0001", "./FCFiles/sample014.fcs", "MEF", 0.85, "sJS1123", "pJS0355, pJS0304, pSC31_3", 187.0829),
    @("S0010", "FC001", "B0001", "./FCFiles/sample015.fcs", "MEF", 0.85, "sJS1123", "pJS0355, pJS0304, pSC31_3", 350),
    @("min",   "FC001", "BMIN",  "./FCFiles/min/sample004.fcs", "MEF", 0.85, "sJS1007", "pJS0143, pJS0130, pSC31_3", 0),
    @("max",   "FC001", "BMAX",  "./FCFiles/max/sample008.fcs", "MEF", 0.85, "sJS1012", "pJS0143, pJS0304, pSC31_3", 0)
)

$r = 2
foreach ($row in $sampleRows) {
    $samples.Cells.Item($r, 1).Value = $row[0]
    $samples.Cells.Item($r, 2).Value = $row[1]
    $samples.Cells.Item($r, 3).Value = $row[2]
    $samples.Cells.Item($r, 4).Value = $row[3]
    $samples.Cells.Item($r, 5).Value = $row[4]
    $samples.Cells.Item($r, 6).Value = $row[5]
    $samples.Cells.Item($r, 7).Value = $row[6]
    $samples.Cells.Item($r, 8).Value = $row[7]
    $samples.Cells.Item($r, 9).Value = $row[8]
    $r++
}
